# Update the dSF (column F) values for the specific rows that changed
# (repull data, push all data, mean calculation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -2
    7  = 3
    8  = -2
    14 = -1
    16 = -3
    18 = -2
    19 = -4
    23 = -4
    36 = -4
    37 = 1
    40 = 0
    41 = -5
    42 = -2
    43 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
